# Task_version_4/src/test/java/data/data.xlsx edit
# - shrink the window height of the workbook view
# - update the stored email address (shared string) in C1
# - move the active selection from C2 to C1
# - widen column C

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shrink the workbook window (best effort - some hosts persist this on the
# window object itself rather than in the saved view XML).
$win = $excel.ActiveWindow
$win.Height = 5088

# The e-mail address text stored in C1 changed.
$ws.Range("C1").Value = "zeinabtest77@mailinator.com"

# The active cell/selection moved from C2 up to C1.
$ws.Range("C1").Select()

# Column C got wider to fit the new text.
$ws.Columns.Item(3).ColumnWidth = 32.6
